# Add two new arrival rows (row 27 and 28) to the "Main Data" sheet,
# mirroring the existing data layout (columns A-L, with K and M left
# present but blank, matching the rest of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: Sunday, Jan 15 - FR2473 from London (STN) - Ryanair B38M (EI-HGW)
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(27, 3).Value = "9:40 AM"
$ws.Cells.Item(27, 4).Value = "FR2473"
$ws.Cells.Item(27, 5).Value = "London"
$ws.Cells.Item(27, 6).Value = "(STN)"
$ws.Cells.Item(27, 7).Value = "Ryanair "
$ws.Cells.Item(27, 8).Value = "B38M"
$ws.Cells.Item(27, 9).Value = "(EI-HGW)"
$ws.Cells.Item(27, 10).Value = "9:25 AM"
$ws.Cells.Item(27, 11).Font.Size = 11
$ws.Cells.Item(27, 12).Value = "0 hours, -15 minutes"
$ws.Cells.Item(27, 13).Font.Size = 11

# Row 28: Sunday, Jan 15 - LO3993 from Warsaw (WAW) - LOT E75S (SP-LIK)
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(28, 3).Value = "2:30 PM"
$ws.Cells.Item(28, 4).Value = "LO3993"
$ws.Cells.Item(28, 5).Value = "Warsaw"
$ws.Cells.Item(28, 6).Value = "(WAW)"
$ws.Cells.Item(28, 7).Value = "LOT "
$ws.Cells.Item(28, 8).Value = "E75S"
$ws.Cells.Item(28, 9).Value = "(SP-LIK)"
$ws.Cells.Item(28, 10).Value = "2:26 PM"
$ws.Cells.Item(28, 11).Font.Size = 11
$ws.Cells.Item(28, 12).Value = "0 hours, -4 minutes"
$ws.Cells.Item(28, 13).Font.Size = 11
